$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows per corrected while-loop check
$ws.Range("A2").Value = 252466
$ws.Range("A3").Value = 252417

# Add new rows 4 and 5
$ws.Range("A4").Value = 252418
$ws.Range("B4").Value = "nessuna compatibilità con alcuna macchina"

$ws.Range("A5").Value = 252980
$ws.Range("B5").Value = "nessuna compatibilità con alcuna macchina"
